# Backlog.xlsx - "Modificacion de ingreso de OC y creacion de reportes de CITI ventas"
#
# Adds two new backlog tasks (rows 114 and 115) to Hoja1, both with
# estado "no comenzado", and moves the active selection to C104
# (matching the post-edit cursor position recorded in the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New backlog rows appended right after the current last row (113).
$ws.Range("A114").Value = "ingreso de mat. Compras pendientes con el mismo id problema para seleccionar"
$ws.Range("B114").Value = "no comenzado"

$ws.Range("A115").Value = "ordenar viajes de hoja de ruta por pendientes y terminados"
$ws.Range("B115").Value = "no comenzado"

# Leave the sheet selection where the author left it when saving.
[void]$ws.Range("C104").Select()
